$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting update: add a no-decimal currency number format and apply
# it to the Total Cost cells (B5, I8, I9) ---
$currencyWholeFormat = '_("$"* #,##0_);_("$"* \(#,##0\);_("$"* "-"??_);_(@_)'

$ws.Range("B5").NumberFormat = $currencyWholeFormat
$ws.Range("I8").NumberFormat = $currencyWholeFormat
$ws.Range("I9").NumberFormat = $currencyWholeFormat

# --- Clear the seeded "0" literal out of the Spots Allocated input cells,
# leaving the cell format (border) intact so the validation below applies
# to a blank, user-entered cell ---
$ws.Range("H8").ClearContents()
$ws.Range("H9").ClearContents()

# --- Row 9's helper columns (L9/M9) were missing; add them to mirror row 8 ---
$ws.Range("L8").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L9").Value = 0

$ws.Range("M8").Copy()
$ws.Range("M9").PasteSpecial(-4122)
$ws.Range("M9").Value = 0

$excel.CutCopyMode = $false

# --- Require a whole number on the Spots Allocated entry column ---
$validationRange = $ws.Range("H8:H9")
$validationRange.Validation.Add(1, 1, 1, "1", "999999999999999000000")
$validationRange.Validation.ErrorTitle = "Error"
$validationRange.Validation.ErrorMessage = "Entry must be a whole number."
$validationRange.Validation.ShowInput = $false
$validationRange.Validation.ShowError = $true
